$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1741.5
$ws.Range("I40").Value = 1266.6666
$ws.Range("J40").Value = 1899.7778
$ws.Range("K40").Value = 1266.6666
$ws.Range("L40").Value = 1899.7778
$ws.Range("M40").Value = -1091.6666
$ws.Range("N40").Value = -2249.7778
$ws.Range("H62").Value = 1923.1
$ws.Range("I62").Value = 2369.8333
$ws.Range("K62").Value = 2369.8333
$ws.Range("M62").Value = -1745.8333
$ws.Range("H64").Value = 2567914.2
$ws.Range("I64").Value = 2962306.8
$ws.Range("J64").Value = 4362.5
$ws.Range("K64").Value = 2962306.8
$ws.Range("L64").Value = 4362.5
$ws.Range("M64").Value = -2962058.8
$ws.Range("N64").Value = -4858.5
$ws.Range("H65").Value = 1923.1
$ws.Range("I65").Value = 2369.8333
$ws.Range("K65").Value = 11849.1665
$ws.Range("M65").Value = -8729.166499999999
$ws.Range("H67").Value = 2567914.2
$ws.Range("I67").Value = 2962306.8
$ws.Range("J67").Value = 4362.5
$ws.Range("K67").Value = 2962306.8
$ws.Range("L67").Value = 4362.5
$ws.Range("M67").Value = -2961448.8
$ws.Range("N67").Value = -6078.5
$ws.Range("H106").Value = 8733.125
$ws.Range("I106").Value = 8902.308000000001
$ws.Range("J106").Value = 8000
$ws.Range("K106").Value = 8902.308000000001
$ws.Range("L106").Value = 8000
$ws.Range("M106").Value = -8271.308000000001
$ws.Range("N106").Value = -9262
$ws.Range("H116").Value = 2392.9333
$ws.Range("I116").Value = 2241.1667
$ws.Range("K116").Value = 2241.1667
$ws.Range("M116").Value = 1200.8333
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 250003700
$ws.Range("I63").Value = 250003700
$ws.Range("K63").Value = 250003700
$ws.Range("M63").Value = -250003014
$ws.Range("H66").Value = 250003700
$ws.Range("I66").Value = 250003700
$ws.Range("K66").Value = 1250018500
$ws.Range("M66").Value = -1250015068
$ws.Range("H117").Value = 49748
$ws.Range("J117").Value = 49748
$ws.Range("L117").Value = 49748
$ws.Range("N117").Value = -58926
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 24390.54
$ws.Range("I26").Value = 8060.625
$ws.Range("J26").Value = 50518.4
$ws.Range("K26").Value = 8060.625
$ws.Range("L26").Value = 50518.4
$ws.Range("M26").Value = -7768.625
$ws.Range("N26").Value = -51102.4
$ws.Range("H86").Value = 25002680
$ws.Range("I86").Value = 2479
$ws.Range("J86").Value = 83336480
$ws.Range("K86").Value = 2479
$ws.Range("L86").Value = 83336480
$ws.Range("M86").Value = -1356
$ws.Range("N86").Value = -83338726
$ws.Range("H89").Value = 25002680
$ws.Range("I89").Value = 2479
$ws.Range("J89").Value = 83336480
$ws.Range("K89").Value = 12395
$ws.Range("L89").Value = 416682400
$ws.Range("M89").Value = -6779
$ws.Range("N89").Value = -416693632
$ws.Range("H105").Value = 4580.4
$ws.Range("J105").Value = 4604.5835
$ws.Range("L105").Value = 4604.5835
$ws.Range("N105").Value = -8098.5835
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2616.3635
$ws.Range("I62").Value = 2480
$ws.Range("J62").Value = 2630
$ws.Range("K62").Value = 2480
$ws.Range("L62").Value = 2630
$ws.Range("M62").Value = -1856
$ws.Range("N62").Value = -3878
$ws.Range("H65").Value = 2616.3635
$ws.Range("I65").Value = 2480
$ws.Range("J65").Value = 2630
$ws.Range("K65").Value = 12400
$ws.Range("L65").Value = 13150
$ws.Range("M65").Value = -9280
$ws.Range("N65").Value = -19390
$ws.Range("H92").Value = 25254.445
$ws.Range("J92").Value = 25254.445
$ws.Range("L92").Value = 25254.445
$ws.Range("N92").Value = -30246.445
$ws.Range("H134").Value = 2987.0833
$ws.Range("I134").Value = 2794.5
$ws.Range("J134").Value = 3950
$ws.Range("K134").Value = 8383.5
$ws.Range("L134").Value = 11850
$ws.Range("M134").Value = -5848.5
$ws.Range("N134").Value = -16920
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 926001.1
$ws.Range("I2").Value = 30
$ws.Range("J2").Value = 1736225.9
$ws.Range("K2").Value = 180
$ws.Range("L2").Value = 10417355.4
$ws.Range("M2").Value = -67
$ws.Range("N2").Value = -10417581.4
$ws.Range("H5").Value = 2597.75
$ws.Range("I5").Value = 274.6
$ws.Range("J5").Value = 4257.143
$ws.Range("K5").Value = 823.8000000000001
$ws.Range("L5").Value = 12771.429
$ws.Range("M5").Value = -711.8000000000001
$ws.Range("N5").Value = -12995.429
$ws.Range("H38").Value = 233.65218
$ws.Range("I38").Value = 230
$ws.Range("J38").Value = 235.6
$ws.Range("K38").Value = 690
$ws.Range("L38").Value = 706.8
$ws.Range("M38").Value = -343
$ws.Range("N38").Value = -1400.8
$ws.Range("H135").Value = 2597.75
$ws.Range("I135").Value = 274.6
$ws.Range("J135").Value = 4257.143
$ws.Range("K135").Value = 2471.4
$ws.Range("L135").Value = 38314.287
$ws.Range("M135").Value = 63.59999999999991
$ws.Range("N135").Value = -43384.287
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20466.416
$ws.Range("I70").Value = 32228.143
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 32228.143
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -31958.143
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 20466.416
$ws.Range("I73").Value = 32228.143
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 32228.143
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -31292.143
$ws.Range("N73").Value = -5872
$ws.Range("H80").Value = 20836144
$ws.Range("I80").Value = 47621348
$ws.Range("J80").Value = 3205.5557
$ws.Range("K80").Value = 47621348
$ws.Range("L80").Value = 3205.5557
$ws.Range("M80").Value = -47620350
$ws.Range("N80").Value = -5201.5557
$ws.Range("H83").Value = 20836144
$ws.Range("I83").Value = 47621348
$ws.Range("J83").Value = 3205.5557
$ws.Range("K83").Value = 238106740
$ws.Range("L83").Value = 16027.7785
$ws.Range("M83").Value = -238101748
$ws.Range("N83").Value = -26011.7785
$ws.Range("H132").Value = 4462.4053
$ws.Range("I132").Value = 4532.6763
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 13598.0289
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -11068.0289
$ws.Range("N132").Value = -16058
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1777.7667
$ws.Range("I16").Value = 1819.0435
$ws.Range("J16").Value = 1642.1428
$ws.Range("K16").Value = 1819.0435
$ws.Range("L16").Value = 1642.1428
$ws.Range("M16").Value = -1649.0435
$ws.Range("N16").Value = -1982.1428
$ws.Range("H122").Value = 10048.75
$ws.Range("I122").Value = 18833.334
$ws.Range("J122").Value = 8021.5386
$ws.Range("K122").Value = 56500.00199999999
$ws.Range("L122").Value = 24064.6158
$ws.Range("M122").Value = -54050.00199999999
$ws.Range("N122").Value = -28964.6158
$ws.Range("H136").Value = 6276.963
$ws.Range("I136").Value = 2736.5833
$ws.Range("K136").Value = 8209.749899999999
$ws.Range("M136").Value = -5659.749899999999
